$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook was re-run to include three additional "spiral" sampling
# schemes (Gaussian-Quadrature moves up next to the other scheme rows, and
# three new Spiral-* rows are inserted after it). The remaining scheme rows
# shift down, and two more rows are appended at the end so every scheme that
# used to be present is still present (now recomputed).
#
# Final row layout (rows 10-19), column order: HKL, then
# [1,1,0] [2,0,0] [2,1,1] [2,2,0] [3,1,0] [2,2,2] [3,2,1] [4,0,0] 2Pairs 4Pairs MaxUnique

$rows = @(
    @{ Row = 10; Label = "Gaussian-Quadrature";          Values = @(1.634508309208395, 2.056121120826005, 1.021648733139336, 1.634508309208395, 0.7277282058053257, 2.018785056086185, 0.7717809265185526, 2.056121120826005, 1.53888492698267, 1.586696618095533, 1.3717620585973) }
    @{ Row = 11; Label = "Spiral-90deg-10rot-5space";    Values = @(1.559426470653133, 3.563089281232255, 0.08727927503991484, 1.559426470653133, 2.022223957573038, 0, 0.3559398217201406, 3.563089281232255, 1.825184278136085, 1.692305374394609, 1.264659801036413) }
    @{ Row = 12; Label = "Spiral-90deg-15rot-5space";    Values = @(1.540973837564633, 3.57031561225003, 0.0874534275581257, 1.540973837564633, 2.027804195942504, 0, 0.3558856388016992, 3.57031561225003, 1.828884519904078, 1.684929178734355, 1.263738785352832) }
    @{ Row = 13; Label = "Spiral-90deg-10rot-3space";    Values = @(1.553384977716487, 3.562501612822962, 0.08738987830544981, 1.553384977716487, 2.023714173948147, 0, 0.3563600782887664, 3.562501612822962, 1.824945745564206, 1.689165361640346, 1.263891786846969) }
    @{ Row = 14; Label = "NoRotation-tilt60deg";         Values = @(4.576016000000003, 2.392376000000007, 0.05374399999999981, 4.576016000000003, 1.071251999999999, 0, 0.3787520000000003, 2.392376000000007, 1.223060000000003, 2.899538000000003, 1.412023333333335) }
    @{ Row = 15; Label = "Rotation-NoTilt";              Values = @(7.650550000000003, 0, 0.04, 7.650550000000003, 0.06, 0, 0.4795874999999998, 0, 0.02, 3.835275000000001, 1.371689583333334) }
    @{ Row = 16; Label = "Rotation-60detTilt";           Values = @(4.725924619264004, 0.4300469861376049, 0.4404363624447993, 4.725924619264004, 0.4607013795840023, 0.424147589324806, 0.7282350641151982, 0.4300469861376049, 0.4352416742912021, 2.580583146777603, 1.201582000145069) }
    @{ Row = 17; Label = "HexGrid-90degTilt5degRes";     Values = @(1.005877367761607, 0.9838600376175911, 0.9959362005835859, 1.005877367761607, 0.9917810608704924, 0.9813296137166525, 0.9965722027548141, 0.9838600376175911, 0.9898981191005884, 0.9978877434310979, 0.9925594138841238) }
    @{ Row = 18; Label = "HexGrid-90degTilt22p5degRes";  Values = @(0.9562114723102951, 1.377496938145186, 0.8806293173617331, 0.9562114723102951, 0.9465545429929211, 1.292504258358423, 0.9513350796612474, 1.377496938145186, 1.12906312775346, 1.042637300031877, 1.067455268138301) }
    @{ Row = 19; Label = "HexGrid-60degTilt5degRes";     Values = @(0.9841436777953232, 1.239648285085894, 0.9527240021869781, 0.9841436777953232, 1.062007812646238, 0.9980678637560858, 0.9582636853054539, 1.239648285085894, 1.096186143636436, 1.040164910715879, 1.032475887795995) }
)

# Rows 17-19 are brand new rows (the sheet used to stop at row 16), so copy
# the "index" cell formatting (bold/centered/bordered) from row 16's column A
# before filling in their values, keeping the same look as the existing rows.
for ($r = 17; $r -le 19; $r++) {
    $ws.Cells.Item(16, 1).Copy($ws.Cells.Item($r, 1))
}

foreach ($rowInfo in $rows) {
    $r = $rowInfo.Row
    $ws.Cells.Item($r, 1).Value = $r - 2
    $ws.Cells.Item($r, 2).Value = $rowInfo.Label

    $col = 3
    foreach ($val in $rowInfo.Values) {
        $ws.Cells.Item($r, $col).Value = $val
        $col = $col + 1
    }
}
